$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 251.35
$ws.Range("I15").Value = 251.35
$ws.Range("K15").Value = 754.05
$ws.Range("M15").Value = -585.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4060.77
$ws.Range("I32").Value = 3824.0938
$ws.Range("J32").Value = 9741
$ws.Range("K32").Value = 3824.0938
$ws.Range("L32").Value = 9741
$ws.Range("M32").Value = -3537.0938
$ws.Range("N32").Value = -10315

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 34000
$ws.Range("J64").Value = 34000
$ws.Range("L64").Value = 34000
$ws.Range("N64").Value = -34496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 34000
$ws.Range("J67").Value = 34000
$ws.Range("L67").Value = 34000
$ws.Range("N67").Value = -35716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 10043741
$ws.Range("I132").Value = 11930216
$ws.Range("J132").Value = 139747.25
$ws.Range("K132").Value = 35790648
$ws.Range("L132").Value = 419241.75
$ws.Range("M132").Value = -35788118
$ws.Range("N132").Value = -424301.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 45000
$ws.Range("I57").Value = 40000
$ws.Range("J57").Value = 50000
$ws.Range("K57").Value = 40000
$ws.Range("L57").Value = 50000
$ws.Range("M57").Value = -39280
$ws.Range("N57").Value = -51440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 35000
$ws.Range("J62").Value = 35000
$ws.Range("L62").Value = 35000
$ws.Range("N62").Value = -36372

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 35000
$ws.Range("J65").Value = 35000
$ws.Range("L65").Value = 105000
$ws.Range("N65").Value = -111864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 48261.668
$ws.Range("J132").Value = 48261.668
$ws.Range("L132").Value = 48261.668
$ws.Range("N132").Value = -58381.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8337848.5
$ws.Range("I134").Value = 4955.2
$ws.Range("J134").Value = 33336528
$ws.Range("K134").Value = 14865.6
$ws.Range("L134").Value = 100009584
$ws.Range("M134").Value = -12330.6
$ws.Range("N134").Value = -100014654

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 52911.8
$ws.Range("J135").Value = 52911.8
$ws.Range("L135").Value = 52911.8
$ws.Range("N135").Value = -63051.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 45000
$ws.Range("I136").Value = 40000
$ws.Range("J136").Value = 50000
$ws.Range("K136").Value = 40000
$ws.Range("L136").Value = 50000
$ws.Range("M136").Value = -34900
$ws.Range("N136").Value = -60200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 15386204
$ws.Range("I58").Value = 27779108
$ws.Range("J58").Value = 1910.7931
$ws.Range("K58").Value = 27779108
$ws.Range("L58").Value = 1910.7931
$ws.Range("M58").Value = -27778905
$ws.Range("N58").Value = -2316.7931

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 32587.412
$ws.Range("J134").Value = 155361.14
$ws.Range("L134").Value = 466083.42
$ws.Range("N134").Value = -471153.42

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 15386204
$ws.Range("I136").Value = 27779108
$ws.Range("J136").Value = 1910.7931
$ws.Range("K136").Value = 83337324
$ws.Range("L136").Value = 5732.379300000001
$ws.Range("M136").Value = -83334774
$ws.Range("N136").Value = -10832.3793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 50629.2
$ws.Range("J140").Value = 50629.2
$ws.Range("L140").Value = 50629.2
$ws.Range("N140").Value = -60989.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 200666.64
$ws.Range("I56").Value = 200666.64
$ws.Range("K56").Value = 200666.64
$ws.Range("M56").Value = -200136.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2352.074
$ws.Range("I122").Value = 1981.4736
$ws.Range("J122").Value = 3232.25
$ws.Range("K122").Value = 5944.4208
$ws.Range("L122").Value = 9696.75
$ws.Range("M122").Value = -3494.4208
$ws.Range("N122").Value = -14596.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 857.56525
$ws.Range("I22").Value = 830
$ws.Range("J22").Value = 887.63635
$ws.Range("K22").Value = 830
$ws.Range("L22").Value = 887.63635
$ws.Range("M22").Value = -535
$ws.Range("N22").Value = -1477.63635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 857.56525
$ws.Range("I27").Value = 830
$ws.Range("J27").Value = 887.63635
$ws.Range("K27").Value = 830
$ws.Range("L27").Value = 887.63635
$ws.Range("M27").Value = -723
$ws.Range("N27").Value = -1101.63635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2562.375
$ws.Range("I61").Value = 2769.3
$ws.Range("J61").Value = 2217.5
$ws.Range("K61").Value = 2769.3
$ws.Range("L61").Value = 2217.5
$ws.Range("M61").Value = -2567.3
$ws.Range("N61").Value = -2621.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1387.1
$ws.Range("I68").Value = 1482.4546
$ws.Range("J68").Value = 1124.875
$ws.Range("K68").Value = 1482.4546
$ws.Range("L68").Value = 1124.875
$ws.Range("M68").Value = -733.4546
$ws.Range("N68").Value = -2622.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1387.1
$ws.Range("I71").Value = 1482.4546
$ws.Range("J71").Value = 1124.875
$ws.Range("K71").Value = 7412.273
$ws.Range("L71").Value = 5624.375
$ws.Range("M71").Value = -3668.273
$ws.Range("N71").Value = -13112.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3002.4
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 3002.4
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 3002.4
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3724.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3002.4
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 3002.4
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 3002.4
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -5498.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1262.0385
$ws.Range("I93").Value = 1300.7222
$ws.Range("J93").Value = 1175
$ws.Range("K93").Value = 1300.7222
$ws.Range("L93").Value = 1175
$ws.Range("M93").Value = -52.72219999999993
$ws.Range("N93").Value = -3671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1510.8214
$ws.Range("I100").Value = 1194.2941
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1194.2941
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -653.2941000000001
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2562.375
$ws.Range("I113").Value = 2769.3
$ws.Range("J113").Value = 2217.5
$ws.Range("K113").Value = 2769.3
$ws.Range("L113").Value = 2217.5
$ws.Range("M113").Value = -599.3000000000002
$ws.Range("N113").Value = -6557.5
